$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the "12 months ended" year headers forward by one year
# (1396/12..1400/12  ->  1397/12..1401/12) across every section header row.
$headerRows = @(8, 16, 26, 35, 43, 52)
foreach ($r in $headerRows) {
    $ws.Range("E$r").Value = "دوازده ماهه منتهی به 1397/12"
    $ws.Range("F$r").Value = "دوازده ماهه منتهی به 1398/12"
    $ws.Range("G$r").Value = "دوازده ماهه منتهی به 1399/12"
    $ws.Range("H$r").Value = "دوازده ماهه منتهی به 1400/12"
    $ws.Range("I$r").Value = "دوازده ماهه منتهی به 1401/12"
}

# --- Update the unit label for "سایر / تخفیفات" under "نرخ فروش" (row 38)
$ws.Range("C38").Value = "/ ریال"

# --- Shift every data row's yearly figures one column to the left and
# append the newly reported value for 1401/12 in column I.
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "-"
$ws.Range("G10").Value = "-"
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0

$ws.Range("E11").Value = 5605950150
$ws.Range("F11").Value = 6185359180
$ws.Range("G11").Value = 6404279210
$ws.Range("H11").Value = 6293613820
$ws.Range("I11").Value = 6258364430

$ws.Range("E12").Value = 5605950150
$ws.Range("F12").Value = 6185359180
$ws.Range("G12").Value = 6404279210
$ws.Range("H12").Value = 6293613820
$ws.Range("I12").Value = 6258364430

$ws.Range("E18").Value = -27000000
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = -2444000
$ws.Range("H18").Value = "-"
$ws.Range("I18").Value = "-"

$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = -7919000
$ws.Range("G19").Value = "-"
$ws.Range("H19").Value = "-"
$ws.Range("I19").Value = "-"

$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0

$ws.Range("E21").Value = 5645000000
$ws.Range("F21").Value = 6316914580
$ws.Range("G21").Value = 6368693180
$ws.Range("H21").Value = 6267074950
$ws.Range("I21").Value = 6233849000

$ws.Range("E22").Value = 5618000000
$ws.Range("F22").Value = 6308995580
$ws.Range("G22").Value = 6366249180
$ws.Range("H22").Value = 6267074950
$ws.Range("I22").Value = 6233849000

$ws.Range("E28").Value = -3014
$ws.Range("F28").Value = -1129
$ws.Range("G28").Value = -391
$ws.Range("H28").Value = "-"
$ws.Range("I28").Value = "-"

$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0

$ws.Range("E30").Value = 654756
$ws.Range("F30").Value = 909038
$ws.Range("G30").Value = 1839441
$ws.Range("H30").Value = 3456701
$ws.Range("I30").Value = 5516295

$ws.Range("E31").Value = 651742
$ws.Range("F31").Value = 907909
$ws.Range("G31").Value = 1839050
$ws.Range("H31").Value = 3456701
$ws.Range("I31").Value = 5516295

$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = "-"
$ws.Range("I37").Value = "-"

$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0

$ws.Range("E39").Value = 116
$ws.Range("F39").Value = 144
$ws.Range("G39").Value = 289
$ws.Range("H39").Value = 552
$ws.Range("I39").Value = 885

$ws.Range("E45").Value = 1796
$ws.Range("F45").Value = 654
$ws.Range("G45").Value = 296
$ws.Range("H45").Value = "-"
$ws.Range("I45").Value = "-"

$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0

$ws.Range("E47").Value = -363165
$ws.Range("F47").Value = -516528
$ws.Range("G47").Value = -773752
$ws.Range("H47").Value = -1620409
$ws.Range("I47").Value = -2070101

$ws.Range("E48").Value = -361369
$ws.Range("F48").Value = -515874
$ws.Range("G48").Value = -773456
$ws.Range("H48").Value = -1620409
$ws.Range("I48").Value = -2070101

$ws.Range("E54").Value = -1218
$ws.Range("F54").Value = -475
$ws.Range("G54").Value = -95
$ws.Range("H54").Value = "-"
$ws.Range("I54").Value = "-"

$ws.Range("E55").Value = "-"
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0

$ws.Range("E56").Value = 291591
$ws.Range("F56").Value = 392510
$ws.Range("G56").Value = 1065689
$ws.Range("H56").Value = 1836292
$ws.Range("I56").Value = 3446194

$ws.Range("E57").Value = 290373
$ws.Range("F57").Value = 392035
$ws.Range("G57").Value = 1065594
$ws.Range("H57").Value = 1836292
$ws.Range("I57").Value = 3446194
